$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pH 10")

# --- H8: value 25 -> 25.9, and pick up the "quote prefix" styled xf (fillId=3, borderId=1, quotePrefix=1) ---
# Excel derives this style when a number is typed with a leading apostrophe into a cell
# that already carries the fillId=3/borderId=1 "data" style. Reproduce the same end state by
# building the quote-prefixed style on a scratch cell (seeded from H8's own current format)
# and pasting that format back onto H8 before writing the real numeric value.
$ws.Range("H8").Value = 25.9
$scratch = $ws.Range("Z1")
$ws.Range("H2").Copy()
$scratch.PasteSpecial(-4122)
$scratch.Value = "'x"
$scratch.Copy()
$ws.Range("H8").PasteSpecial(-4122)
$scratch.Clear()

# --- fill in previously empty SiO2 [mg/l] (H) measurements ---
$ws.Range("H9").Value = 55.7
$ws.Range("H10").Value = 27.3
$ws.Range("H11").Value = 63.2
$ws.Range("H12").Value = 29.6
$ws.Range("H13").Value = 71
$ws.Range("H14").Value = 30.7

# --- re-key the time points: 6 -> 5.5 (rows 11-12), 7 -> 6 (rows 13-14) ---
$ws.Range("A11").Value = 5.5
$ws.Range("A12").Value = 5.5
$ws.Range("A13").Value = 6
$ws.Range("A14").Value = 6

# --- new rows 15 & 16: a new "6.5" timepoint, Feed + Permeate ---
$ws.Range("A14:I14").Copy()
$ws.Range("A15:I15").PasteSpecial(-4122)
$ws.Range("A14:I14").Copy()
$ws.Range("A16:I16").PasteSpecial(-4122)

# A15/A16 hold "6.5" as literal text (matches the rest of column A's look while staying
# distinct from the numeric 6 above it) - build it as a formula-derived string, then
# freeze it down to a plain value so it lands as a shared string, not a live formula.
$ws.Range("A15").Formula = '="6.5"'
$ws.Range("A15").Copy()
$ws.Range("A15").PasteSpecial(-4163)
$ws.Range("A16").Formula = '="6.5"'
$ws.Range("A16").Copy()
$ws.Range("A16").PasteSpecial(-4163)

$ws.Range("B15").Value = "Feed"
$ws.Range("C15").Formula = '=CONCATENATE("M","-9.2-",A15,"-F")'
$ws.Range("H15").Value = 78.5

$ws.Range("B16").Value = "Permeate"
$ws.Range("C16").Formula = '=CONCATENATE("M","-9.2-",A16,"-P")'
$ws.Range("H16").Value = 32.4

# --- view/selection state: "pH 10" becomes the active/selected tab ---
$ws.Activate()
$ws.Range("H17").Select()

Write-Host "done"
